$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 0.1740869811147974
$ws.Cells.Item(2, 5).Value = 0.1584551540435015
$ws.Cells.Item(2, 6).Value = 1.586058537684792
$ws.Cells.Item(2, 7).Value = 0.002415559989363013
$ws.Cells.Item(2, 10).Value = 0.1913305446091726
$ws.Cells.Item(2, 15).Value = 3.975293305671983
$ws.Cells.Item(3, 4).Value = 0.1738283514954659
$ws.Cells.Item(3, 5).Value = 0.1560976871574269
$ws.Cells.Item(3, 6).Value = 1.537274724656413
$ws.Cells.Item(3, 7).Value = 0.002420193372918766
$ws.Cells.Item(3, 10).Value = 0.1859881629286377
$ws.Cells.Item(3, 15).Value = 3.829925653828354
$ws.Cells.Item(4, 4).Value = 0.1737352727056418
$ws.Cells.Item(4, 5).Value = 0.1547365595057606
$ws.Cells.Item(4, 6).Value = 1.508292599236015
$ws.Cells.Item(4, 7).Value = 0.002423190390455585
$ws.Cells.Item(4, 10).Value = 0.1828259976755433
$ws.Cells.Item(4, 15).Value = 3.743161388050737
$ws.Cells.Item(5, 4).Value = 0.173713882169821
$ws.Cells.Item(5, 5).Value = 0.154203567062531
$ws.Cells.Item(5, 6).Value = 1.496725533069409
$ws.Cells.Item(5, 7).Value = 0.002424450076757445
$ws.Cells.Item(5, 10).Value = 0.181566964411239
$ws.Cells.Item(5, 15).Value = 3.70842818185713
$ws.Cells.Item(6, 4).Value = 0.1737113296677322
$ws.Cells.Item(6, 5).Value = 0.1541163722461789
$ws.Cells.Item(6, 6).Value = 1.494819505439409
$ws.Cells.Item(6, 7).Value = 0.002424661568398314
$ws.Cells.Item(6, 10).Value = 0.1813596857493067
$ws.Cells.Item(6, 15).Value = 3.702698365185711
$ws.Cells.Item(7, 4).Value = 0.1737349172369989
$ws.Cells.Item(7, 5).Value = 0.1547292836603447
$ws.Cells.Item(7, 6).Value = 1.50813561742865
$ws.Cells.Item(7, 7).Value = 0.002423207223660427
$ws.Cells.Item(7, 10).Value = 0.1828088983185054
$ws.Cells.Item(7, 15).Value = 3.742690441195975
$ws.Cells.Item(8, 4).Value = 0.1739841729762546
$ws.Cells.Item(8, 5).Value = 0.1576243489919804
$ws.Cells.Item(8, 6).Value = 1.569035736272326
$ws.Cells.Item(8, 7).Value = 0.002417126089511913
$ws.Cells.Item(8, 10).Value = 0.1894639112098275
$ws.Cells.Item(8, 15).Value = 3.92465161631327
$ws.Cells.Item(9, 4).Value = 0.1749940350984289
$ws.Cells.Item(9, 5).Value = 0.1639891454149165
$ws.Cells.Item(9, 6).Value = 1.696215314610924
$ws.Cells.Item(9, 7).Value = 0.002406401938534081
$ws.Cells.Item(9, 10).Value = 0.2034575054671848
$ws.Cells.Item(9, 15).Value = 4.301396925444863
$ws.Cells.Item(10, 4).Value = 0.176053429039797
$ws.Cells.Item(10, 5).Value = 0.1690884747248944
$ws.Cells.Item(10, 6).Value = 1.794459758914087
$ws.Cells.Item(10, 7).Value = 0.00239924673843737
$ws.Cells.Item(10, 10).Value = 0.2143235631690743
$ws.Cells.Item(10, 15).Value = 4.590573542432764
$ws.Cells.Item(11, 4).Value = 0.1766042903178189
$ws.Cells.Item(11, 5).Value = 0.1715010782485393
$ws.Cells.Item(11, 6).Value = 1.840215008519493
$ws.Cells.Item(11, 7).Value = 0.00239614705872108
$ws.Cells.Item(11, 10).Value = 0.2193960977187288
$ws.Cells.Item(11, 15).Value = 4.724870106639401
$ws.Cells.Item(12, 4).Value = 0.1768227941905849
$ws.Cells.Item(12, 5).Value = 0.1724280870052368
$ws.Cells.Item(12, 6).Value = 1.857695457956879
$ws.Cells.Item(12, 7).Value = 0.002394995481506987
$ws.Cells.Item(12, 10).Value = 0.2213357164881842
$ws.Cells.Item(12, 15).Value = 4.776123705914642
$ws.Cells.Item(13, 4).Value = 0.176775295088845
$ws.Cells.Item(13, 5).Value = 0.1722278419399643
$ws.Cells.Item(13, 6).Value = 1.853923867927165
$ws.Cells.Item(13, 7).Value = 0.002395242508834146
$ws.Cells.Item(13, 10).Value = 0.2209171483286099
$ws.Cells.Item(13, 15).Value = 4.765067556405256
$ws.Cells.Item(14, 4).Value = 0.1766220683053774
$ws.Cells.Item(14, 5).Value = 0.1715770748055547
$ws.Cells.Item(14, 6).Value = 1.841650044738742
$ws.Cells.Item(14, 7).Value = 0.002396051873487624
$ws.Cells.Item(14, 10).Value = 0.2195552945536008
$ws.Cells.Item(14, 15).Value = 4.729078768254283
$ws.Cells.Item(15, 4).Value = 0.1765295021121389
$ws.Cells.Item(15, 5).Value = 0.1711802087939205
$ws.Cells.Item(15, 6).Value = 1.83415205034737
$ws.Cells.Item(15, 7).Value = 0.002396550520847136
$ws.Cells.Item(15, 10).Value = 0.2187235669911161
$ws.Cells.Item(15, 15).Value = 4.707086574762002
$ws.Cells.Item(16, 4).Value = 0.1760188154244986
$ws.Cells.Item(16, 5).Value = 0.1689326789280017
$ws.Cells.Item(16, 6).Value = 1.791491030024503
$ws.Cells.Item(16, 7).Value = 0.002399452423562754
$ws.Cells.Item(16, 10).Value = 0.2139946788175706
$ws.Cells.Item(16, 15).Value = 4.581852524519547
$ws.Cells.Item(17, 4).Value = 0.1757231760738733
$ws.Cells.Item(17, 5).Value = 0.1675777197095343
$ws.Cells.Item(17, 6).Value = 1.765592947740515
$ws.Cells.Item(17, 7).Value = 0.002401272324906944
$ws.Cells.Item(17, 10).Value = 0.2111269263639315
$ws.Cells.Item(17, 15).Value = 4.505731717013873
$ws.Cells.Item(18, 4).Value = 0.1755596213571522
$ws.Cells.Item(18, 5).Value = 0.1668071234688853
$ws.Cells.Item(18, 6).Value = 1.750797082196556
$ws.Cells.Item(18, 7).Value = 0.00240233370465939
$ws.Cells.Item(18, 10).Value = 0.2094896528481058
$ws.Cells.Item(18, 15).Value = 4.462207568522899
$ws.Cells.Item(19, 4).Value = 0.1755053593350553
$ws.Cells.Item(19, 5).Value = 0.1665477124916208
$ws.Cells.Item(19, 6).Value = 1.745804605181746
$ws.Cells.Item(19, 7).Value = 0.002402695584620141
$ws.Cells.Item(19, 10).Value = 0.2089373882817682
$ws.Cells.Item(19, 15).Value = 4.447515346699447
$ws.Cells.Item(20, 4).Value = 0.1757539758419568
$ws.Cells.Item(20, 5).Value = 0.1677210524654029
$ws.Cells.Item(20, 6).Value = 1.768339483154421
$ws.Cells.Item(20, 7).Value = 0.002401077080862745
$ws.Cells.Item(20, 10).Value = 0.2114309417491569
$ws.Cells.Item(20, 15).Value = 4.513808127811672
$ws.Cells.Item(21, 4).Value = 0.1766668059815757
$ws.Cells.Item(21, 5).Value = 0.1717678565881826
$ws.Cells.Item(21, 6).Value = 1.845250978317722
$ws.Cells.Item(21, 7).Value = 0.002395813541768803
$ws.Cells.Item(21, 10).Value = 0.2199547937790811
$ws.Cells.Item(21, 15).Value = 4.739638709769963
$ws.Cells.Item(22, 4).Value = 0.1773211216157335
$ws.Cells.Item(22, 5).Value = 0.1744908541517773
$ws.Cells.Item(22, 6).Value = 1.896414790630274
$ws.Cells.Item(22, 7).Value = 0.002392502880711286
$ws.Cells.Item(22, 10).Value = 0.2256350364452402
$ws.Cells.Item(22, 15).Value = 4.889555756342702
$ws.Cells.Item(23, 4).Value = 0.176966621572177
$ws.Cells.Item(23, 5).Value = 0.1730303690042234
$ws.Cells.Item(23, 6).Value = 1.869025219674029
$ws.Cells.Item(23, 7).Value = 0.002394258045954197
$ws.Cells.Item(23, 10).Value = 0.2225933271111273
$ws.Cells.Item(23, 15).Value = 4.809328534291183
$ws.Cells.Item(24, 4).Value = 0.1757400312938202
$ws.Cells.Item(24, 5).Value = 0.1676562256042757
$ws.Cells.Item(24, 6).Value = 1.767097483987953
$ws.Cells.Item(24, 7).Value = 0.002401165303850142
$ws.Cells.Item(24, 10).Value = 0.211293460785356
$ws.Cells.Item(24, 15).Value = 4.510156039134927
$ws.Cells.Item(25, 4).Value = 0.1746650798688378
$ws.Cells.Item(25, 5).Value = 0.1621932848350518
$ws.Cells.Item(25, 6).Value = 1.660971659900724
$ws.Cells.Item(25, 7).Value = 0.002409175396274585
$ws.Cells.Item(25, 10).Value = 0.1995698711110521
$ws.Cells.Item(25, 15).Value = 4.197320322167002
